$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns.
# Values that parse as plain numbers are given a leading apostrophe so
# Excel stores them as text, preserving the source's exact decimal
# formatting (e.g. trailing zeros like "1.000" or "0.9980").

$ws.Range("D2").Value = '29.333.56'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '1.847.50'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("D4").Value = '''0.9988'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''240.06'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").Value = '''0.6268'
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = '''0.9984'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '''0.07603'
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").Value = '''0.2909'
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("D10").Value = '''24.69'
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").Value = '''0.07736'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = '''5.018'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '''0.6787'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = '''0.00001054'
$ws.Range("E14").Value = '  -3.40%  '
$ws.Range("D15").Value = '''83.07'
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").Value = '''6.132'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = '29.370.64'
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("D18").Value = '''229.67'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").Value = '''12.34'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").Value = '''0.9987'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '''7.486'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '''158.35'
$ws.Range("E23").Value = '  +0.80%  '
$ws.Range("D24").Value = '''0.1385'
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("D25").Value = '''8.432'
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("D26").Value = '''17.68'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '''1.432'
$ws.Range("E27").Value = '  +8.91%  '
$ws.Range("D28").Value = '''1.467'
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").Value = '''0.05606'
$ws.Range("E29").Value = '  -2.07%  '
$ws.Range("D30").Value = '''4.103'
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("D31").Value = '''4.067'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").Value = '''1.161'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = '''1.831'
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").Value = '''0.6944'
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("D35").Value = '''2.580'
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("D36").Value = '1.234.30'
$ws.Range("E36").Value = '  +1.39%  '
$ws.Range("D37").Value = '''0.01797'
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").Value = '''2.727'
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("D39").Value = '''6.367'
$ws.Range("E39").Value = '  -1.79%  '
$ws.Range("D40").Value = '''0.9018'
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").Value = '''0.9980'
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").Value = '''101.34'
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("D43").Value = '''65.50'
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("D44").Value = '''7.180'
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").Value = '''0.00000000117'
$ws.Range("E45").Value = '  -2.74%  '
$ws.Range("D46").Value = '''0.3990'
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("D49").Value = '''0.1145'
$ws.Range("E49").Value = '  +1.41%  '
$ws.Range("D50").Value = '''0.05695'
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").Value = '''0.4621'
$ws.Range("E51").Value = '  -0.14%  '

# Rows 47/48: RenderToken and EnergySwap swap positions, with updated price/volume
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''9.024'
$ws.Range("E47").Value = '  +0.39%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.686'
$ws.Range("E48").Value = '  +0.17%  '
